$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.001.79'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').Value = '2.414.70'
$ws.Range('E3').Value = '  -0.23%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '552.55'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('E6').Value = '  -1.01%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  +3.79%  '
$ws.Range('E9').Value = '  -1.76%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '5.68'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.55%  '
$ws.Range('E11').Value = '  -0.89%  '
$ws.Range('E12').Value = '  -1.62%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '25.33'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +2.32%  '
$ws.Range('D14').Value = '2.845.42'
$ws.Range('E14').Value = '  -0.12%  '
$ws.Range('D15').Value = '59.918.89'
$ws.Range('E15').Value = '  +0.20%  '
$ws.Range('E16').Value = '  -1.35%  '
$ws.Range('D17').Value = '2.409.68'
$ws.Range('E17').Value = '  +0.80%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.32'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.90%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.42'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '328.23'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.94%  '
$ws.Range('E21').Value = '  -3.16%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '65.89'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.11%  '
$ws.Range('E24').Value = '  +3.13%  '
$ws.Range('E25').Value = '  +0.72%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').Value = '  +0.51%  '
$ws.Range('D28').Value = '0.0₃0775'
$ws.Range('E28').Value = '  -1.66%  '
$ws.Range('E29').Value = '  -2.13%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '169.02'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.08%  '
$ws.Range('E31').Value = '  -3.80%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '18.62'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.55%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('E36').Value = '  +0.09%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.19'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.07%  '
$ws.Range('E38').Value = '  -2.01%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '322.17'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.72%  '
$ws.Range('E40').Value = '  -4.47%  '
$ws.Range('E41').Value = '  -1.69%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '140.55'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.55%  '
$ws.Range('E43').Value = '  +0.60%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '19.58'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.11%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0515'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.59%  '
$ws.Range('E46').Value = '  +1.09%  '
$ws.Range('E47').Value = '  -1.28%  '
$ws.Range('E48').Value = '  -7.00%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '11.04'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('E50').Value = '  -3.32%  '
